# "Added size on disk for each model"
#
# The underlying data/formatting for the "Size on disk (MB)" column was
# already present in the workbook; what this revision captures is simply
# the author re-saving the file after reviewing it, which Excel records as
# a moved selection/active cell on the "Kitti" sheet (from F20 to O13).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kitti")

$ws.Activate() | Out-Null
$ws.Range("O13").Select() | Out-Null
